# Daily attendance processing - 2025-10-09 11:41:47
# Reorders the "Recorded By" (column G) value for specific rows so that
# "System" appears after the other recorder instead of before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(3,4,6,11,12,13,14,15,30,33,38,39,40,41,42,57,58,60,65,66,67,68,69,86,89,93,95,112,115,119,121,138,141,145,147)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $value = $cell.Value2
    $parts = $value -split ',\s*'
    if ($parts.Count -eq 2 -and $parts[0] -eq 'System') {
        $cell.Value = "$($parts[1]), $($parts[0])"
    }
}
